# Add two derived percentage columns (I: homeless-population %, J: total-population %)
# to the COVID-19 homeless impact dataset, mirroring the upstream commit
# "Added percentage of population affected".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
# Column I: percent of the homeless population (cumulative_homeless_cases / 151278)
# Column J: percent of the total population (total_cases / 39510000)
# NOTE: write J before I so the shared-strings table picks up the same
# ordering as the source workbook (total_population_by_percentage first,
# then total_homless_population_by_percentage).
$ws.Range("J1").Value = "total_population_by_percentage "
$ws.Range("I1").Value = "total_homless_population_by_percentage "

# --- First data row: literal (non-shared) formulas ---
$ws.Range("I2").Formula = "=(H2/151278)*100"
$ws.Range("J2").Formula = "=(G2/39510000)*100"

# --- Bulk-fill the rest in two batches (mirrors the source file's shared-formula groups) ---
$ws.Range("I3:I66").Formula = "=(H3/151278)*100"
$ws.Range("J3:J66").Formula = "=(G3/39510000)*100"

$ws.Range("I67:I91").Formula = "=(H67/151278)*100"
$ws.Range("J67:J90").Formula = "=(G67/39510000)*100"
$ws.Range("J91").Formula = "=(G91/39510000)*100"

# --- Column widths for the two new columns (closest reachable to the source's 33.22 / 25.55 char widths) ---
$ws.Columns.Item(9).ColumnWidth = 32.3
$ws.Columns.Item(10).ColumnWidth = 24.6

# --- Match the saved selection/active cell from the source edit ---
$ws.Range("I23").Select()
